$wb = $excel.ActiveWorkbook

# Sheet "展览": F4 1336 -> 1346, F5 647 -> 651
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1346
$ws1.Range("F5").Value = 651

# Sheet "全部类型": F4 1336 -> 1346, F6 647 -> 651
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1346
$ws4.Range("F6").Value = 651
